# Sync attendance_reports: normalize "Recorded By" (column G) ordering so
# that entries reorder with "System" moved to the front of the
# comma-separated list (equivalent to reversing the list, since "System"
# is always the last entry in the original data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notlike "*,*") { continue }

    $rawParts = $val.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    if ($parts -notcontains "System") { continue }

    $n = $parts.Count
    $reversed = @()
    for ($i = $n - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }

    $newVal = [string]::Join(", ", $reversed)
    $cell.Value2 = $newVal
}
